$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "Documentatie"
$ws.Range("C9").Value = "Implementeren"
$ws.Range("C10").Value = "Puntejes op de I"

$ws.Range("C10").Select()
